# Remove both comments attached to "John" (commentRangeStart/End +
# commentReference runs in word/document.xml, and the comment bodies in
# word/comments.xml). Deleting every Comment object removes the comment
# markup from the main story as well as the comment part's content.
$d = $word.ActiveDocument

while ($d.Comments.Count -gt 0) {
    $d.Comments(1).Delete()
}

# The removal of the comment anchors is the "last edit" Word records, so
# the hidden _GoBack bookmark follows it: it moves from the trailing empty
# paragraph at the end of the document to the point right after "John"
# (immediately before " Doe"). Re-adding a bookmark with the same name
# moves it (Word keeps only one bookmark per name).
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Execute("John", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$goBack = $d.Range($rng.End, $rng.End)
$d.Bookmarks.Add("_GoBack", $goBack) | Out-Null
